$wb = $excel.ActiveWorkbook

# The workbook used to ship two near-identical "Users" template sheets
# (xr1 / xr2). Consolidate on a single, properly named sheet.
$wb.Worksheets.Item("xr2").Delete() | Out-Null

$ws = $wb.Worksheets.Item("xr1")
$ws.Name = "Users"

# New "Kelas" (class) column next to Password.
$ws.Range("E1").Value = "Kelas"
$ws.Range("E1").AddComment("Optional")

# The old per-column explanation comments for Role/Password are no longer
# needed on this simplified template.
$ws.Range("C1").Comment.Delete()
$ws.Range("D1").Comment.Delete()

# Refresh the sample row with generic placeholder data instead of a real
# student's info.
$ws.Hyperlinks.Delete()
$ws.Range("A2").Value = "test@gmail.com"
$ws.Range("B2").Value = "Testing"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:test@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("A16").Select() | Out-Null

Write-Output "done"
